$d = $word.ActiveDocument

# The document title starts with a run containing "GIT " followed by a
# run containing "REBASE". Replace the text of just that first run
# (the first 4 characters of the document) so that the title reads
# " REBASE" instead of "GIT REBASE" (the leading "GIT" is removed but
# the trailing space of that run is kept), without touching the
# separate "REBASE" run that follows it.
$r = $d.Range(0, 4)
$r.Text = " "
